# "add more items in elimination list and to add list"
#
# The "HypothyroidismAdd" sheet (the active sheet) currently has its
# recipe rows laid out with a blank row between every entry (rows
# 4, 6, 8, ... 28) - left over from a previous edit to the elimination
# sheet. Tidy it up by removing those blank spacer rows so the recipe
# rows become contiguous (rows 3-15), matching the already-compacted
# "HypothyroidismEliminate" sheet.
#
# Deleting from the bottom-most blank row upward means each deletion's
# upward shift never disturbs the row numbers of blank rows still
# waiting to be removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 27; $r -ge 3; $r -= 2) {
    $ws.Rows($r).Delete()
}
